# BOT; UPDATE DATA
# Appends the 2020-05-13 ("43964") consultation-count row to the
# "相談件数" sheet, pushing the trailing footnote row down by one, and
# keeps the print area / used-range / selection in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Make room for the new data row just above the footnote row (old row 109)
# -- this shifts the footnote row (and its style/shared-string content)
# down to row 110 automatically.
$ws.Rows.Item(109).Insert()

# New day's figures.
$ws.Range("A109").Value = 43964
$ws.Range("B109").Value = 297
$ws.Range("C109").Value = 36845
$ws.Range("D109").Value = 51
$ws.Range("E109").Value = 7488

# Extend the print area to cover the newly-added row.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "=相談件数!`$A`$1:`$E`$111"
    }
}

# Match the saved selection on the data sheet.
$ws.Range("A109").Select()

Write-Host "done"
